# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker/period detail table (rows 16-37, columns C:E) is re-sorted so
# that all workers for period "1710" come first (rows 16-26), followed by
# the same workers for period "1711" (rows 27-37), keeping the original
# per-worker ordering within each period block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Period 1710 block (rows 16-26)
$ws.Range("C16").Value = "1100334373"
$ws.Range("D16").Value = "EDUARDO SANEZ ANGEL"
$ws.Range("E16").Value = "1710"

$ws.Range("C17").Value = "73190652"
$ws.Range("D17").Value = "JOSE LUIS BLANCO CABARCAS"
$ws.Range("E17").Value = "1710"

$ws.Range("C18").Value = "1047432610"
$ws.Range("D18").Value = "EDWIN RAFAEL MUÑOZ DE ANGEL"
$ws.Range("E18").Value = "1710"

$ws.Range("C19").Value = "1052981975"
$ws.Range("D19").Value = "FRANKY ALI RIOS MUÑOZ"
$ws.Range("E19").Value = "1710"

$ws.Range("C20").Value = "1052992439"
$ws.Range("D20").Value = "MIGUEL ANGEL JIMENEZ MARQUEZ"
$ws.Range("E20").Value = "1710"

$ws.Range("C21").Value = "1052968659"
$ws.Range("D21").Value = "DIONISIO JOSE RIOS MUÑOZ"
$ws.Range("E21").Value = "1710"

$ws.Range("C22").Value = "1052960660"
$ws.Range("D22").Value = "RAFAEL OVIDIO RIOS MUÑOZ"
$ws.Range("E22").Value = "1710"

$ws.Range("C23").Value = "1052946538"
$ws.Range("D23").Value = "EVERTO MANUEL DE LAS OSSA ATENCIA"
$ws.Range("E23").Value = "1710"

$ws.Range("C24").Value = "1101445961"
$ws.Range("D24").Value = "PEDRO LUIS TORRES ZUÑIGA"
$ws.Range("E24").Value = "1710"

$ws.Range("C25").Value = "1052956555"
$ws.Range("D25").Value = "JUAN FRANCISCO MEZA LOPEZ"
$ws.Range("E25").Value = "1710"

$ws.Range("C26").Value = "1103220109"
$ws.Range("D26").Value = "DIEGO ARMANDO WILCHEZ NOVOA"
$ws.Range("E26").Value = "1710"

# Period 1711 block (rows 27-37)
$ws.Range("C27").Value = "1100334373"
$ws.Range("D27").Value = "EDUARDO SANEZ ANGEL"
$ws.Range("E27").Value = "1711"

$ws.Range("C28").Value = "73190652"
$ws.Range("D28").Value = "JOSE LUIS BLANCO CABARCAS"
$ws.Range("E28").Value = "1711"

$ws.Range("C29").Value = "1047432610"
$ws.Range("D29").Value = "EDWIN RAFAEL MUÑOZ DE ANGEL"
$ws.Range("E29").Value = "1711"

$ws.Range("C30").Value = "1052981975"
$ws.Range("D30").Value = "FRANKY ALI RIOS MUÑOZ"
$ws.Range("E30").Value = "1711"

$ws.Range("C31").Value = "1052992439"
$ws.Range("D31").Value = "MIGUEL ANGEL JIMENEZ MARQUEZ"
$ws.Range("E31").Value = "1711"

$ws.Range("C32").Value = "1052968659"
$ws.Range("D32").Value = "DIONISIO JOSE RIOS MUÑOZ"
$ws.Range("E32").Value = "1711"

$ws.Range("C33").Value = "1052960660"
$ws.Range("D33").Value = "RAFAEL OVIDIO RIOS MUÑOZ"
$ws.Range("E33").Value = "1711"

$ws.Range("C34").Value = "1052946538"
$ws.Range("D34").Value = "EVERTO MANUEL DE LAS OSSA ATENCIA"
$ws.Range("E34").Value = "1711"

$ws.Range("C35").Value = "1101445961"
$ws.Range("D35").Value = "PEDRO LUIS TORRES ZUÑIGA"
$ws.Range("E35").Value = "1711"

$ws.Range("C36").Value = "1052956555"
$ws.Range("D36").Value = "JUAN FRANCISCO MEZA LOPEZ"
$ws.Range("E36").Value = "1711"

$ws.Range("C37").Value = "1103220109"
$ws.Range("D37").Value = "DIEGO ARMANDO WILCHEZ NOVOA"
$ws.Range("E37").Value = "1711"
